$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 40: B40 was stored as text "4"; fix it to a real number 4.
$ws.Range("B40").Value = 4

# Append new row 41 with the new annotation data.
$ws.Range("A41").Value = "Ying Tang"
$ws.Range("C41").Value = "无"
$ws.Range("D41").Value = "DIS"
$ws.Range("E41").Value = "MET"
$ws.Range("F41").Value = "df7b0ece-3727-4ec6-95ce-2a2839e398ed"
$ws.Range("G41").Value = "SkhQHMW0W_annotated.xlsx"
$ws.Range("H41").Value = "This is necessary to get an idea of total amount of communication that was sufficient to reach perplexity 72.24 at the end of 40-th epoch."

# B41 needs to stay as TEXT "3" (not be auto-coerced to a number like a plain
# Value assignment would do). Build it as a text formula result in a scratch
# cell, then paste-special as values so the destination cell keeps a text
# type without picking up a stray number format / style.
$ws.Range("Z1").Formula = "=""3"""
$ws.Range("Z1").Copy()
$ws.Range("B41").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
